$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "72.478.09"
$ws.Range("E2").Value = "  +4.39%  "

# Row 3
$ws.Range("D3").Value = "4.039.32"
$ws.Range("E3").Value = "  +3.44%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.67"
$ws.Range("E5").Value = "  -2.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.68"
$ws.Range("E6").Value = "  +1.42%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.724"
$ws.Range("E7").Value = "  +18.32%  "

# Row 8
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.759"
$ws.Range("E9").Value = "  +5.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  +1.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000325"
$ws.Range("E11").Value = "  -2.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.97"
$ws.Range("E12").Value = "  +11.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.97"
$ws.Range("E13").Value = "  +6.94%  "

# Row 14
$ws.Range("D14").Value = "4.690.54"
$ws.Range("E14").Value = "  +3.51%  "

# Row 15
$ws.Range("D15").Value = "4.041.92"
$ws.Range("E15").Value = "  +3.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.10"
$ws.Range("E16").Value = "  +6.87%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.17"
$ws.Range("E17").Value = "  +0.89%  "

# Row 18
$ws.Range("E18").Value = "  -0.10%  "

# Row 19
$ws.Range("E19").Value = "  -1.73%  "

# Row 20
$ws.Range("D20").Value = "72.357.46"
$ws.Range("E20").Value = "  +4.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.11"
$ws.Range("E21").Value = "  +2.77%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "104.00"
$ws.Range("E22").Value = "  +17.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.60"
$ws.Range("E23").Value = "  +6.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.57"
$ws.Range("E24").Value = "  +2.23%  "

# Row 25
$ws.Range("E25").Value = "  -0.81%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.53"
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.08"
$ws.Range("E27").Value = "  +4.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.06"
$ws.Range("E28").Value = "  +4.61%  "

# Row 29
$ws.Range("E29").Value = "  +2.10%  "

# Row 30
$ws.Range("E30").Value = "  +10.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.67"
$ws.Range("E31").Value = "  +3.91%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("E32").Value = "  +1.11%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "668.63"
$ws.Range("E33").Value = "  -3.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.80"
$ws.Range("E34").Value = "  +14.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.69"
$ws.Range("E35").Value = "  +0.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.19"
$ws.Range("E36").Value = "  +5.78%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.434"
$ws.Range("E37").Value = "  -1.58%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0862"
$ws.Range("E38").Value = "  +1.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.54"
$ws.Range("E39").Value = "  +13.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.151"
$ws.Range("E40").Value = "  +0.27%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0493"
$ws.Range("E42").Value = "  +2.41%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.20"
$ws.Range("E44").Value = "  +2.90%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.157"
$ws.Range("E45").Value = "  +12.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  -2.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.50"
$ws.Range("E47").Value = "  +4.26%  "

# Row 48
$ws.Range("E48").Value = "  +2.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +6.75%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.10"
$ws.Range("E50").Value = "  +1.63%  "

# Row 51
$ws.Range("E51").Value = "  +1.02%  "
